$wb = $excel.ActiveWorkbook

# --- Overview sheet ---
# "Latest HO Xliff Generate Date" column (G) for the two rows whose
# timestamp was refreshed by the new handback run.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G2").Value = "2016-09-05 16:19:24"
$wsOverview.Range("G5").Value = "2016-09-05 16:19:24"

# --- zh-cn sheet ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
# Priority changed from human translation (ht) to machine translation (mt)
$wsZhCn.Range("E2").Value = "mt"
$wsZhCn.Range("E5").Value = "mt"
# Correspond Handoff Datetime refreshed
$wsZhCn.Range("H2").Value = "2016-09-05 16:19:19"
$wsZhCn.Range("H5").Value = "2016-09-05 16:19:19"
# Correspond Handback DateTime refreshed
$wsZhCn.Range("K2").Value = "2016-09-05 16:19:36"
$wsZhCn.Range("K5").Value = "2016-09-05 16:19:36"

# --- de-de sheet ---
$wsDeDe = $wb.Worksheets.Item("de-de")
# Priority changed from human translation (ht) to machine translation (mt)
$wsDeDe.Range("E2").Value = "mt"
$wsDeDe.Range("E5").Value = "mt"
# Correspond Handoff Datetime refreshed (shares the same timestamp as the
# Overview sheet's "Latest HO Xliff Generate Date")
$wsDeDe.Range("H2").Value = "2016-09-05 16:19:24"
$wsDeDe.Range("H5").Value = "2016-09-05 16:19:24"
# Correspond Handback DateTime refreshed
$wsDeDe.Range("K2").Value = "2016-09-05 16:19:44"
$wsDeDe.Range("K5").Value = "2016-09-05 16:19:44"
